# Updated cryptos list on Thu Apr 11 05:40:19 UTC 2024 with GitHub Actions
#
# This script refreshes the Price (D) and Volume(1h) (E) columns of the
# crypto table, and fixes up the two rows (18 & 19) whose coin/link/price
# values were previously swapped (WrappedBTC <-> WrappedEther).
#
# Several "Price" values (column D) are plain decimals such as "611.96"
# or "1.00". If assigned directly via .Value, Excel COM auto-detects them
# as numbers, which both changes their cell type and silently drops
# formatting such as the trailing zero in "1.00" (-> 1) or "17.40" (-> 17.4).
# To preserve the original text semantics we prefix those values with a
# leading apostrophe (Excel's classic "force text" marker) and then reset
# the cell Style afterwards so the visual/style-level "quote prefix"
# indicator that Excel adds doesn't leak into the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Value,
        [bool]$ForceText = $false
    )

    $range = $ws.Range($Cell)
    if ($ForceText) {
        $range.Value = "'" + $Value
        $range.Style = "Normal"
    } else {
        $range.Value = $Value
    }
}

# --- Row 2 : Bitcoin ---
Set-TextValue "D2" "70.578.57"
Set-TextValue "E2" "  +1.77%  "

# --- Row 3 : Ethereum ---
Set-TextValue "D3" "3.563.60"
Set-TextValue "E3" "  +0.83%  "

# --- Row 4 : TetherUSD ---
Set-TextValue "E4" "  +0.00%  "

# --- Row 5 : BNB ---
Set-TextValue "D5" "611.96" $true
Set-TextValue "E5" "  +5.35%  "

# --- Row 6 : Solana ---
Set-TextValue "D6" "173.21" $true
Set-TextValue "E6" "  +0.72%  "

# --- Row 7 ---
Set-TextValue "E7" "  +1.34%  "

# --- Row 8 ---
Set-TextValue "D8" "3.559.90"
Set-TextValue "E8" "  +0.84%  "

# --- Row 9 ---
Set-TextValue "E9" "  -0.01%  "

# --- Row 10 ---
Set-TextValue "E10" "  +3.31%  "

# --- Row 11 ---
Set-TextValue "D11" "7.48" $true
Set-TextValue "E11" "  +12.43%  "

# --- Row 12 ---
Set-TextValue "D12" "0.587" $true
Set-TextValue "E12" "  -0.13%  "

# --- Row 13 ---
Set-TextValue "D13" "46.67" $true
Set-TextValue "E13" "  -1.67%  "

# --- Row 14 ---
Set-TextValue "E14" "  +0.76%  "

# --- Row 15 ---
Set-TextValue "D15" "4.141.01"
Set-TextValue "E15" "  +1.15%  "

# --- Row 16 ---
Set-TextValue "E16" "  -2.42%  "

# --- Row 17 : BitcoinCash ---
Set-TextValue "D17" "615.90" $true
Set-TextValue "E17" "  -2.11%  "

# --- Row 18 & 19 : WrappedBTC / WrappedEther were swapped; fix order + data ---
Set-TextValue "B18" "WrappedEther"
Set-TextValue "C18" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D18" "3.560.98"
Set-TextValue "E18" "  +0.77%  "

Set-TextValue "B19" "WrappedBTC"
Set-TextValue "C19" "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue "D19" "70.671.88"
Set-TextValue "E19" "  +1.98%  "

# --- Row 20 ---
Set-TextValue "E20" "  -2.01%  "

# --- Row 21 ---
Set-TextValue "D21" "17.40" $true
Set-TextValue "E21" "  -0.75%  "

# --- Row 22 ---
Set-TextValue "D22" "0.884" $true
Set-TextValue "E22" "  -0.74%  "

# --- Row 23 ---
Set-TextValue "D23" "9.38" $true
Set-TextValue "E23" "  -16.46%  "

# --- Row 24 ---
Set-TextValue "D24" "16.10" $true
Set-TextValue "E24" "  +0.39%  "

# --- Row 25 ---
Set-TextValue "D25" "97.14" $true
Set-TextValue "E25" "  -0.79%  "

# --- Row 26 ---
Set-TextValue "E26" "  +0.57%  "

# --- Row 27 ---
Set-TextValue "D27" "1.00" $true
Set-TextValue "E27" "  -0.14%  "

# --- Row 28 ---
Set-TextValue "D28" "2.63" $true
Set-TextValue "E28" "  -0.88%  "

# --- Row 29 ---
Set-TextValue "E29" "  +1.59%  "

# --- Row 30 ---
Set-TextValue "D30" "9.07" $true
Set-TextValue "E30" "  -2.56%  "

# --- Row 31 ---
Set-TextValue "D31" "8.52" $true
Set-TextValue "E31" "  -0.83%  "

# --- Row 32 ---
Set-TextValue "E32" "  -3.53%  "

# --- Row 33 ---
Set-TextValue "D33" "6.98" $true
Set-TextValue "E33" "  -0.72%  "

# --- Row 34 ---
Set-TextValue "E34" "  -2.49%  "

# --- Row 35 ---
Set-TextValue "D35" "576.37" $true
Set-TextValue "E35" "  -8.77%  "

# --- Row 36 ---
Set-TextValue "D36" "3.69" $true
Set-TextValue "E36" "  +5.36%  "

# --- Row 37 ---
Set-TextValue "E37" "  -2.06%  "

# --- Row 38 ---
Set-TextValue "E38" "  +0.21%  "

# --- Row 39 ---
Set-TextValue "D39" "0.0479" $true
Set-TextValue "E39" "  +4.91%  "

# --- Row 40 ---
Set-TextValue "D40" "57.31" $true
Set-TextValue "E40" "  -0.14%  "

# --- Row 41 ---
Set-TextValue "E41" "  +0.06%  "

# --- Row 42 ---
Set-TextValue "E42" "  +3.74%  "

# --- Row 43 ---
Set-TextValue "D43" "3.392.00"
Set-TextValue "E43" "  -0.20%  "

# --- Row 44 ---
Set-TextValue "E44" "  -2.97%  "

# --- Row 45 ---
Set-TextValue "D45" "33.19" $true
Set-TextValue "E45" "  +0.41%  "

# --- Row 46 ---
Set-TextValue "D46" "2.98" $true
Set-TextValue "E46" "  +7.34%  "

# --- Row 47 ---
Set-TextValue "E47" "  +0.81%  "

# --- Row 48 ---
Set-TextValue "E48" "  +1.52%  "

# --- Row 49 ---
Set-TextValue "E49" "  -0.15%  "

# --- Row 50 ---
Set-TextValue "D50" "133.68" $true
Set-TextValue "E50" "  +1.62%  "
